$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14-24 down to 15-25 (bottom-up so sources aren't clobbered before
# being read), copying both the functionality text (column B) and the status
# cell formatting (column C fill/border).
for ($r = 24; $r -ge 14; $r--) {
    $src = $r
    $dst = $r + 1

    $ws.Range("B$dst").Value = $ws.Range("B$src").Value()

    $ws.Range("C$src").Copy()
    $ws.Range("C$dst").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Row 25 is brand new, give it the same row height as the rest of the table.
$ws.Rows.Item(25).RowHeight = 15

# "Admin Ban user when on profile" (row 13) is now finished: mark it solved
# (green), matching the style already used for other solved rows.
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Row 14 becomes the newly finished "Admin Remove Ban" functionality, also
# solved.
$ws.Range("B14").Value = "Admin Remove Ban"
$ws.Range("C2").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep column A a clean sequential list (1..24) after the shift.
for ($r = 14; $r -le 25; $r++) {
    $ws.Range("A$r").Value = $r - 1
}

$ws.Range("F12").Select()
